# Microsite Education Script completed
# Appends new interview-history rows to AMSIN, BETA and AMS sheets, and
# corrects an existing row (AMS!27) whose timestamp/style were slightly off.

function Add-HistoryRow {
    param($ws, $r, $dateText, $timestamp, $cycleText, $d, $e, $f, $g)

    # Column A: plain text date label (quote-prefix keeps it General/text,
    # not auto-parsed into a date).
    $ws.Cells.Item($r, 1).Value = "'" + $dateText

    # Column B: the numeric Excel serial timestamp, shown with the sheet's
    # custom date/time number format.
    $cb = $ws.Cells.Item($r, 2)
    $cb.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cb.Value = $timestamp

    # Column C: text identifier/tag for the run.
    $ws.Cells.Item($r, 3).Value = "'" + $cycleText

    # Columns D-G: numeric metrics.
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AMSIN: add rows 52-56
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Add-HistoryRow $wsAmsin 52 "2022-08-02" 44775.65561364583 "165_fstcycle"  165 165 0 4.41
Add-HistoryRow $wsAmsin 53 "2022-08-03" 44776.66674120371 "165_scndcycle" 165 165 0 33.28
Add-HistoryRow $wsAmsin 54 "2022-08-04" 44777.38928859954 "165_finalrun"  165 165 0 4.49
Add-HistoryRow $wsAmsin 55 "2022-08-22" 44795.66193466435 "166fstcycle"   165 165 0 4.59
Add-HistoryRow $wsAmsin 56 "2022-08-23" 44796.90335474537 "166cyclescnd"  165 165 0 4.47

# ---------------------------------------------------------------------
# BETA: add rows 25-26
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
Add-HistoryRow $wsBeta 25 "2022-08-04" 44777.55889136574 "165beta" 165 165 0 4.7
Add-HistoryRow $wsBeta 26 "2022-08-24" 44797.53743711806 "166_beta" 165 165 0 4.8

# ---------------------------------------------------------------------
# AMS: fix row 27 (timestamp correction + explicit "Normal" styling that
# matches the other data rows), then add rows 28-29
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Re-assert the text values (quote-prefixed, like every other text cell in
# this sheet) and restyle the numeric cells to "Normal" so the whole row
# carries an explicit style instead of the sheet's implicit default.
$wsAms.Cells.Item(27, 1).Value = "'2022-07-14"
$wsAms.Cells.Item(27, 3).Value = "'164_live"
$wsAms.Range("D27:G27").Style = "Normal"
$wsAms.Cells.Item(27, 2).Value = 44756.81726386574

Add-HistoryRow $wsAms 28 "2022-08-04" 44777.81277537037 "165_live" 165 165 0 4.84
Add-HistoryRow $wsAms 29 "2022-08-24" 44797.91952252982 "166_live" 165 165 0 4.69
